# covariants_test.xlsx - "changed photographs after review"
# Adds 16 new "Mestioni" specimen photographs to Sheet1, extends the
# tribe/subtribe classification formulas to recognise them, and drops in
# a small Gel/Mes/Ech legend next to the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. New rows 32-47: additional specimen photo IDs (column A)
# ---------------------------------------------------------------------
$newIds = @(
    "apopra_gelina_mesgab_g.JPG",
    "apshir_gelina_mesagr_g.JPG",
    "artfor_gelina_messph_g.JPG",
    "atetsi_gelina_mesate_g.JPG",
    "bicver_gelina_meslym_g.JPG",
    "coeful_gelina_mescoe_g.JPG",
    "crespi_gelina_mescer_g.JPG",
    "cryban_gelina_mesmes_g.JPG",
    "dicsp_gelina_mesglo_g.JPG",
    "melsp_gelina_mesmel_g.JPG",
    "necsp_gelina_mesgor_g.JPG",
    "odospi_gelina_mesisc_g.JPG",
    "parsap_gelina_mesbar_g.JPG",
    "picmel_gelina_mesnem_g.JPG",
    "steaxi_gelina_mesvag_g.JPG",
    "steins_gelina_mesnem_g.JPG"
)

$startRow = 32
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newIds[$i]
}
$endRow = $startRow + $newIds.Length - 1   # 47

# ---------------------------------------------------------------------
# 2. TRIBE column (B): recognise "mes" -> "Mestioni" as well as "ech"
# ---------------------------------------------------------------------
# Existing rows 2-31 keep using the same formula range, just with the
# updated logic (adds the nested Mestioni test; behaviour for existing
# rows is unchanged since none of them contain "mes").
$ws.Range("B2:B31").Formula = '=IF(ISNUMBER(SEARCH("ech",A2)), "Echthrini", IF(ISNUMBER(SEARCH("mes",A2)), "Mestioni", "Gelini"))'

# New rows: row 32 on its own, then the rest as a block.
$ws.Range("B32").Formula = '=IF(ISNUMBER(SEARCH("ech",A32)), "Echthrini", IF(ISNUMBER(SEARCH("mes",A32)), "Mestioni", "Gelini"))'
$ws.Range("B33:B" + $endRow).Formula = '=IF(ISNUMBER(SEARCH("ech",A33)), "Echthrini", IF(ISNUMBER(SEARCH("mes",A33)), "Mestioni", "Gelini"))'

# ---------------------------------------------------------------------
# 3. OPEN/CLOSED column (D): extend the existing logic down to row 47
# ---------------------------------------------------------------------
$ws.Range("D32:D" + $endRow).Formula = '=IF(ISNUMBER(SEARCH("_O",A32)), "OPEN", "CLOSED")'

# ---------------------------------------------------------------------
# 4. Small legend table next to the header (F8:G10)
# ---------------------------------------------------------------------
$ws.Range("F8").Value = "Gel"
$ws.Range("G8").Value = "Gelini"
$ws.Range("F9").Value = "Mes"
$ws.Range("G9").Value = "Mesostini"
$ws.Range("F10").Value = "Ech"
$ws.Range("G10").Value = "Echthrini"

# ---------------------------------------------------------------------
# 5. Selection cursor moves to D2 (matches the saved view state)
# ---------------------------------------------------------------------
$ws.Range("D2").Select()
